$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.999.06'
$ws.Range("E2").Value = '  -0.18%  '
$ws.Range("D3").Value = '1.912.87'
$ws.Range("E3").Value = '  +0.62%  '
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").Value = "'324.62"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.82%  '
$ws.Range("E6").Value = '  -0.34%  '
$ws.Range("D7").Value = "'0.4598"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.43%  '
$ws.Range("D8").Value = "'0.3869"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.71%  '
$ws.Range("D9").Value = "'0.07822"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.01%  '
$ws.Range("D10").Value = "'0.9887"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.59%  '
$ws.Range("D11").Value = "'21.88"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.54%  '
$ws.Range("D12").Value = '1.884.37'
$ws.Range("E12").Value = '  +0.76%  '
$ws.Range("D13").Value = "'5.772"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.64%  '
$ws.Range("D14").Value = "'7.010"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.81%  '
$ws.Range("D15").Value = "'0.07048"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.22%  '
$ws.Range("D16").Value = "'87.33"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.31%  '
$ws.Range("E17").Value = '  -0.19%  '
$ws.Range("D18").Value = "'0.000009907"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.23%  '
$ws.Range("D19").Value = "'17.03"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("D20").Value = "'0.9999"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.40%  '
$ws.Range("D21").Value = '29.009.45'
$ws.Range("E21").Value = '  -0.15%  '
$ws.Range("D22").Value = "'5.366"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.64%  '
$ws.Range("D23").Value = "'11.09"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.51%  '
$ws.Range("D24").Value = '2.137.76'
$ws.Range("E24").Value = '  +1.15%  '
$ws.Range("D25").Value = "'2.087"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.04%  '
$ws.Range("D26").Value = "'156.11"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.20%  '
$ws.Range("D27").Value = "'19.37"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.15%  '
$ws.Range("D28").Value = "'5.832"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -2.98%  '
$ws.Range("D29").Value = "'118.30"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.21%  '
$ws.Range("D30").Value = "'1.854"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.97%  '
$ws.Range("D31").Value = "'0.09308"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.51%  '
$ws.Range("D32").Value = "'0.8778"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.13%  '
$ws.Range("D33").Value = "'5.182"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.82%  '
$ws.Range("D34").Value = "'1.309"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.89%  '
$ws.Range("D35").Value = "'3.129"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -4.83%  '
$ws.Range("D36").Value = "'0.05755"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.88%  '
$ws.Range("D37").Value = "'1.166"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.24%  '
$ws.Range("D38").Value = "'0.02084"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.46%  '
$ws.Range("D39").Value = "'0.9989"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.39%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = "'7.641"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.80%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = "'0.5665"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.87%  '
$ws.Range("D42").Value = "'0.1806"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +1.19%  '
$ws.Range("D43").Value = "'0.000003027"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +89.90%  '
$ws.Range("E44").Value = '  -2.10%  '
$ws.Range("E45").Value = '  -1.08%  '
$ws.Range("E46").Value = '  -3.42%  '
$ws.Range("D47").Value = "'0.5300"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.45%  '
$ws.Range("D48").Value = "'0.06926"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.74%  '
$ws.Range("B49").Value = 'MXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D49").Value = "'2.567"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.10%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = "'1.832"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.61%  '
$ws.Range("D51").Value = "'112.50"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.36%  '
